$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1387.6875
$ws.Range("I33").Value = 923.2308
$ws.Range("K33").Value = 923.2308
$ws.Range("M33").Value = -694.2308

$ws.Range("H100").Value = 1799.3
$ws.Range("I100").Value = 1602.3846
$ws.Range("J100").Value = 2165
$ws.Range("K100").Value = 1602.3846
$ws.Range("L100").Value = 2165
$ws.Range("M100").Value = -1061.3846
$ws.Range("N100").Value = -3247

$ws.Range("H115").Value = 1900.9
$ws.Range("I115").Value = 1686.4286
$ws.Range("K115").Value = 5059.2858
$ws.Range("M115").Value = -3492.2858

$ws.Range("H129").Value = 1095.4546
$ws.Range("I129").Value = 453.17648
$ws.Range("J129").Value = 1499.8518
$ws.Range("K129").Value = 1359.52944
$ws.Range("L129").Value = 4499.555399999999
$ws.Range("M129").Value = 3640.47056
$ws.Range("N129").Value = -14499.5554

$ws.Range("H138").Value = 2045.0807
$ws.Range("I138").Value = 1596.0714
$ws.Range("J138").Value = 2414.853
$ws.Range("K138").Value = 4788.2142
$ws.Range("L138").Value = 7244.559
$ws.Range("M138").Value = 351.7857999999997
$ws.Range("N138").Value = -17524.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 441360.2
$ws.Range("I32").Value = 499804.84
$ws.Range("K32").Value = 499804.84
$ws.Range("M32").Value = -499517.84

$ws.Range("H61").Value = 2249.0667
$ws.Range("I61").Value = 1979.6923
$ws.Range("K61").Value = 1979.6923
$ws.Range("M61").Value = -1767.6923

$ws.Range("H102").Value = 2909
$ws.Range("I102").Value = 2998.889
$ws.Range("K102").Value = 2998.889
$ws.Range("M102").Value = -1376.889

$ws.Range("H110").Value = 1062.2
$ws.Range("I110").Value = 932.6667
$ws.Range("J110").Value = 1256.5
$ws.Range("K110").Value = 932.6667
$ws.Range("L110").Value = 1256.5
$ws.Range("M110").Value = 1112.3333
$ws.Range("N110").Value = -5346.5

$ws.Range("H136").Value = 2249.0667
$ws.Range("I136").Value = 1979.6923
$ws.Range("K136").Value = 5939.0769
$ws.Range("M136").Value = -3389.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5684123
$ws.Range("I105").Value = 5954729
$ws.Range("K105").Value = 5954729
$ws.Range("M105").Value = -5952982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5952.849
$ws.Range("I31").Value = 1180.5
$ws.Range("J31").Value = 14492.842
$ws.Range("K31").Value = 1180.5
$ws.Range("L31").Value = 14492.842
$ws.Range("M31").Value = -885.5
$ws.Range("N31").Value = -15082.842

$ws.Range("H34").Value = 5952.849
$ws.Range("I34").Value = 1180.5
$ws.Range("J34").Value = 14492.842
$ws.Range("K34").Value = 1180.5
$ws.Range("L34").Value = 14492.842
$ws.Range("M34").Value = -978.5
$ws.Range("N34").Value = -14896.842

$ws.Range("H58").Value = 1434.238
$ws.Range("J58").Value = 1744.5
$ws.Range("L58").Value = 1744.5
$ws.Range("N58").Value = -2150.5

$ws.Range("H105").Value = 1508
$ws.Range("I105").Value = 1185.7142
$ws.Range("K105").Value = 1185.7142
$ws.Range("M105").Value = 561.2858000000001

$ws.Range("H136").Value = 1434.238
$ws.Range("J136").Value = 1744.5
$ws.Range("L136").Value = 5233.5
$ws.Range("N136").Value = -10333.5

$ws.Range("H141").Value = 216666.67
$ws.Range("J141").Value = 220000
$ws.Range("L141").Value = 220000
$ws.Range("N141").Value = -230360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11000
$ws.Range("I3").Value = 11000
$ws.Range("K3").Value = 33000
$ws.Range("M3").Value = -32888

$ws.Range("H11").Value = 823.75
$ws.Range("I11").Value = 931.6667
$ws.Range("K11").Value = 2795.0001
$ws.Range("M11").Value = -2655.0001

$ws.Range("H21").Value = 2999.5
$ws.Range("I21").Value = 2999
$ws.Range("K21").Value = 8997
$ws.Range("M21").Value = -8824

$ws.Range("H36").Value = 1250
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -1331
$ws.Range("N36").Value = -6338

$ws.Range("H129").Value = 1851.85
$ws.Range("I129").Value = 480
$ws.Range("J129").Value = 2093.9412
$ws.Range("K129").Value = 1440
$ws.Range("L129").Value = 6281.823600000001
$ws.Range("M129").Value = 3560
$ws.Range("N129").Value = -16281.8236

$ws.Range("H131").Value = 1046.6976
$ws.Range("I131").Value = 434
$ws.Range("J131").Value = 1127.3158
$ws.Range("K131").Value = 1302
$ws.Range("L131").Value = 3381.9474
$ws.Range("M131").Value = 3738
$ws.Range("N131").Value = -13461.9474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H120").Value = 32700
$ws.Range("I120").Value = 30000
$ws.Range("J120").Value = 34050
$ws.Range("K120").Value = 30000
$ws.Range("L120").Value = 34050
$ws.Range("M120").Value = -25162
$ws.Range("N120").Value = -43726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5576.815
$ws.Range("I22").Value = 1390
$ws.Range("J22").Value = 8039.647
$ws.Range("K22").Value = 1390
$ws.Range("L22").Value = 8039.647
$ws.Range("M22").Value = -1095
$ws.Range("N22").Value = -8629.647000000001

$ws.Range("H27").Value = 5576.815
$ws.Range("I27").Value = 1390
$ws.Range("J27").Value = 8039.647
$ws.Range("K27").Value = 1390
$ws.Range("L27").Value = 8039.647
$ws.Range("M27").Value = -1283
$ws.Range("N27").Value = -8253.647000000001

$ws.Range("H68").Value = 1604.2916
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2334.3333
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2334.3333
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3832.3333

$ws.Range("H71").Value = 1604.2916
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2334.3333
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 11671.6665
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -19159.6665

$ws.Range("H100").Value = 2294.8823
$ws.Range("I100").Value = 2064.125
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2064.125
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1523.125
$ws.Range("N100").Value = -3582

$ws.Range("H136").Value = 13891315
$ws.Range("J136").Value = 27780346
$ws.Range("L136").Value = 83341038
$ws.Range("N136").Value = -83346138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 42856.5
$ws.Range("J137").Value = 42856.5
$ws.Range("L137").Value = 42856.5
$ws.Range("N137").Value = -53056.5

$ws.Range("H139").Value = 87665
$ws.Range("J139").Value = 87665
$ws.Range("L139").Value = 87665
$ws.Range("N139").Value = -97945
